$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "330.85"; E = "0.67%"; G = "20" }
    @{ Row = 3; D = "41.21"; E = "2.14%"; G = "20" }
    @{ Row = 4; D = "5.686"; E = "-2.77%"; G = "20" }
    @{ Row = 5; D = "0.08039"; E = "-1.11%"; G = "20" }
    @{ Row = 6; D = "2.023"; E = "2.88%"; G = "20" }
    @{ Row = 7; D = "8.728"; E = "-0.38%"; G = "20" }
    @{ Row = 8; D = "4.536"; E = "-1.52%"; G = "20" }
    @{ Row = 9; D = "2.978"; E = "1.23%"; G = "20" }
    @{ Row = 10; D = "0.9228"; E = "-2.44%"; G = "20" }
    @{ Row = 11; D = "0.1258"; E = "-5.25%"; G = "20" }
    @{ Row = 12; D = "0.1946"; E = "-2.63%"; G = "20" }
    @{ Row = 13; D = ""; E = "-3.19%"; G = "20" }
    @{ Row = 14; D = "0.09341"; E = "0.10%"; G = "20" }
    @{ Row = 15; D = "0.03759"; E = "9.02%"; G = "20" }
    @{ Row = 16; D = "0.1052"; E = "9.36%"; G = "20" }
    @{ Row = 17; D = "0.001304"; E = "-0.67%"; G = "20" }
    @{ Row = 18; D = "0.006225"; E = "-1.36%"; G = "20" }
    @{ Row = 19; D = "3.366"; E = "0.03%"; G = "20" }
    @{ Row = 20; D = "0.3479"; E = "-1.70%"; G = "20" }
    @{ Row = 21; D = "0.1421"; E = "1.26%"; G = "20" }
    @{ Row = 22; D = "0.2661"; E = "10.19%"; G = "20" }
    @{ Row = 23; D = "0.04441"; E = "0.41%"; G = "20" }
    @{ Row = 24; D = "0.001263"; E = "-0.01%"; G = "20" }
    @{ Row = 25; D = "0.004286"; E = "-2.62%"; G = "20" }
    @{ Row = 26; D = "0.0001245"; E = "13.83%"; G = "20" }
    @{ Row = 27; D = ""; E = ""; G = "20" }
    @{ Row = 28; D = ""; E = ""; G = "20" }
    @{ Row = 29; D = ""; E = ""; G = "20" }
    @{ Row = 30; D = ""; E = ""; G = "20" }
    @{ Row = 31; D = ""; E = ""; G = "20" }
    @{ Row = 32; D = ""; E = ""; G = "20" }
    @{ Row = 33; D = ""; E = ""; G = "20" }
    @{ Row = 34; D = ""; E = ""; G = "20" }
    @{ Row = 35; D = ""; E = ""; G = "20" }
    @{ Row = 36; D = ""; E = ""; G = "20" }
    @{ Row = 37; D = ""; E = ""; G = "20" }
    @{ Row = 38; D = ""; E = ""; G = "20" }
    @{ Row = 39; D = "0.02855"; E = "15.05%"; G = "20" }
    @{ Row = 40; D = "0.05478"; E = "3.65%"; G = "20" }
    @{ Row = 41; D = "0.007794"; E = "3.94%"; G = "20" }
    @{ Row = 42; D = "0.009968"; E = "11.22%"; G = "20" }
    @{ Row = 43; D = "0.1413"; E = "-1.53%"; G = "20" }
    @{ Row = 44; D = "0.002129"; E = "3.50%"; G = "20" }
    @{ Row = 45; D = "0.01188"; E = "12.66%"; G = "20" }
    @{ Row = 46; D = "0.00006787"; E = "-1.80%"; G = "20" }
    @{ Row = 47; D = "0.00000000753"; E = "0.20%"; G = "20" }
    @{ Row = 48; D = "0.002289"; E = "34.35%"; G = "20" }
    @{ Row = 49; D = "0.003019"; E = "-13.91%"; G = "20" }
    @{ Row = 50; D = "0.00002109"; E = "0.20%"; G = "20" }
    @{ Row = 51; D = "0.0002008"; E = "0.20%"; G = "20" }
)


foreach ($u in $updates) {
    $row = $u.Row
    if ($u.D -ne "") {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.E -ne "") {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
    if ($u.G -ne "") {
        $cell = $ws.Cells.Item($row, 7)
        $cell.NumberFormat = "@"
        $cell.Value = $u.G
    }
}
